$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.023388977412426
$ws.Range("D2").Value = 1.027504567682766
$ws.Range("E2").Value = 1.027012581899041
$ws.Range("F2").Value = 1.033593909381782
$ws.Range("I2").Value = 1.029967066348059
$ws.Range("J2").Value = 1.028569472579588
$ws.Range("K2").Value = 1.030324102740825
$ws.Range("L2").Value = 1.029833551620001
$ws.Range("M2").Value = 1.036395810392662
$ws.Range("N2").Value = 1.030030159209221
# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.024444389513681
$ws.Range("D3").Value = 1.028253779698649
$ws.Range("E3").Value = 1.028015397125965
$ws.Range("F3").Value = 1.034840654544225
$ws.Range("I3").Value = 1.03014789802563
$ws.Range("J3").Value = 1.029262793063152
$ws.Range("K3").Value = 1.030881212111393
$ws.Range("L3").Value = 1.030643474224827
$ws.Range("M3").Value = 1.037450400348963
$ws.Range("N3").Value = 1.030724464287395
# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.025127421028426
$ws.Range("D4").Value = 1.028738459555429
$ws.Range("E4").Value = 1.028664763832025
$ws.Range("F4").Value = 1.03564761550342
$ws.Range("I4").Value = 1.030263505897426
$ws.Range("J4").Value = 1.029710986511492
$ws.Range("K4").Value = 1.031240927721233
$ws.Range("L4").Value = 1.031167421229479
$ws.Range("M4").Value = 1.038132471672365
$ws.Range("N4").Value = 1.03117329422184
# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.025414594063063
$ws.Range("D5").Value = 1.02894219199714
$ws.Range("E5").Value = 1.028937871981171
$ws.Range("F5").Value = 1.035986918940797
$ws.Range("I5").Value = 1.030311771569476
$ws.Range("J5").Value = 1.02989930375838
$ws.Range("K5").Value = 1.031391967241002
$ws.Range("L5").Value = 1.031387658069773
$ws.Range("M5").Value = 1.038419138987156
$ws.Range("N5").Value = 1.031361878900818
# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.025462813250379
$ws.Range("D6").Value = 1.028976397973393
$ws.Range("E6").Value = 1.028983734746031
$ws.Range("F6").Value = 1.036043892860082
$ws.Range("I6").Value = 1.030319855887539
$ws.Range("J6").Value = 1.029930916995497
$ws.Range("K6").Value = 1.031417316589842
$ws.Range("L6").Value = 1.031424635015294
$ws.Range("M6").Value = 1.038467267290624
$ws.Range("N6").Value = 1.031393537032358
# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.025131258143577
$ws.Range("D7").Value = 1.028741181944941
$ws.Range("E7").Value = 1.028668412667141
$ws.Range("F7").Value = 1.035652149066939
$ws.Range("I7").Value = 1.030264152145666
$ws.Range("J7").Value = 1.029713503221621
$ws.Range("K7").Value = 1.031242946645622
$ws.Range("L7").Value = 1.031170364164759
$ws.Range("M7").Value = 1.038136302430507
$ws.Range("N7").Value = 1.031175814505985
# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.023745636489954
$ws.Range("D8").Value = 1.027757789932337
$ws.Range("E8").Value = 1.027351388929019
$ws.Range("F8").Value = 1.034015205100289
$ws.Range("I8").Value = 1.030028469522719
$ws.Range("J8").Value = 1.028803872760153
$ws.Range("K8").Value = 1.030512540038106
$ws.Range("L8").Value = 1.030107295034644
$ws.Range("M8").Value = 1.036752280561881
$ws.Range("N8").Value = 1.030264892264926
# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.021304812270692
$ws.Range("D9").Value = 1.026024105126535
$ws.Range("E9").Value = 1.025034290934735
$ws.Range("F9").Value = 1.031132436757547
$ws.Range("I9").Value = 1.029602433447484
$ws.Range("J9").Value = 1.027197685668222
$ws.Range("K9").Value = 1.02921957270282
$ws.Range("L9").Value = 1.028233060170746
$ws.Range("M9").Value = 1.034310975343086
$ws.Range("N9").Value = 1.028656424203092
# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.019678116835516
$ws.Range("D10").Value = 1.024867791543917
$ws.Range("E10").Value = 1.023492025832287
$ws.Range("F10").Value = 1.029211683046513
$ws.Range("I10").Value = 1.029311201301024
$ws.Range("J10").Value = 1.026124673086638
$ws.Range("K10").Value = 1.02835364146181
$ws.Range("L10").Value = 1.026982917367737
$ws.Range("M10").Value = 1.032681712106467
$ws.Range("N10").Value = 1.027581887820567
# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.01897385584341
$ws.Range("D11").Value = 1.024366976281895
$ws.Range("E11").Value = 1.022824792598642
$ws.Range("F11").Value = 1.028380216850588
$ws.Range("I11").Value = 1.029183386312148
$ws.Range("J11").Value = 1.025659518926282
$ws.Range("K11").Value = 1.0279777479994
$ws.Range("L11").Value = 1.02644143603553
$ws.Range("M11").Value = 1.031975799986176
$ws.Range("N11").Value = 1.027116073087949
# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.018712277565689
$ws.Range("D12").Value = 1.024180932952327
$ws.Range("E12").Value = 1.022577039112523
$ws.Range("F12").Value = 1.028071406645019
$ws.Range("I12").Value = 1.029135653352078
$ws.Range("J12").Value = 1.02548665975569
$ws.Range("K12").Value = 1.027837983350003
$ws.Range("L12").Value = 1.026240281361514
$ws.Range("M12").Value = 1.03171352685107
$ws.Range("N12").Value = 1.026942968437507
# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.01876838629967
$ws.Range("D13").Value = 1.024220840697753
$ws.Range("E13").Value = 1.022630179157889
$ws.Range("F13").Value = 1.028137645998775
$ws.Range("I13").Value = 1.029145903853054
$ws.Range("J13").Value = 1.02552374226693
$ws.Range("K13").Value = 1.027867969724816
$ws.Range("L13").Value = 1.026283430828751
$ws.Range("M13").Value = 1.031769788333499
$ws.Range("N13").Value = 1.026980103610166
# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.018952233387191
$ws.Range("D14").Value = 1.024351598252282
$ws.Range("E14").Value = 1.022804311451537
$ws.Range("F14").Value = 1.028354689831486
$ws.Range("I14").Value = 1.029179445928722
$ws.Range("J14").Value = 1.025645231970523
$ws.Range("K14").Value = 1.027966197887479
$ws.Range("L14").Value = 1.02642480902613
$ws.Range("M14").Value = 1.03195412175964
$ws.Range("N14").Value = 1.027101765843073
# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.019065509692043
$ws.Range("D15").Value = 1.024432159890713
$ws.Range("E15").Value = 1.022911611581527
$ws.Range("F15").Value = 1.028488422072315
$ws.Range("I15").Value = 1.029200078284167
$ws.Range("J15").Value = 1.02572007515323
$ws.Range("K15").Value = 1.02802670080797
$ws.Range("L15").Value = 1.026511913561504
$ws.Range("M15").Value = 1.03206768689084
$ws.Range("N15").Value = 1.027176715311684
# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.01972485854829
$ws.Range("D16").Value = 1.024901026425291
$ws.Range("E16").Value = 1.02353632006371
$ws.Range("F16").Value = 1.029266869486846
$ws.Range("I16").Value = 1.029319647947632
$ws.Range("J16").Value = 1.026155532598817
$ws.Range("K16").Value = 1.028378568484464
$ws.Range("L16").Value = 1.027018850288448
$ws.Range("M16").Value = 1.032728551991825
$ws.Range("N16").Value = 1.027612791156792
# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.02013847928738
$ws.Range("D17").Value = 1.025195101115681
$ws.Range("E17").Value = 1.023928337758565
$ws.Range("F17").Value = 1.029755230177121
$ws.Range("I17").Value = 1.029394193185112
$ws.Range("J17").Value = 1.026428540820849
$ws.Range("K17").Value = 1.028599034296607
$ws.Range("L17").Value = 1.027336795012463
$ws.Range("M17").Value = 1.03314297907833
$ws.Range("N17").Value = 1.027886187081814
# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.020379747820424
$ws.Range("D18").Value = 1.025366617933248
$ws.Range("E18").Value = 1.024157050800744
$ws.Range("F18").Value = 1.030040105279333
$ws.Range("I18").Value = 1.02943750921049
$ws.Range("J18").Value = 1.026587730518904
$ws.Range("K18").Value = 1.028727537620023
$ws.Range("L18").Value = 1.02752223124422
$ws.Range("M18").Value = 1.033384665946188
$ws.Range("N18").Value = 1.028045602847499
# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.020462015968007
$ws.Range("D19").Value = 1.025425098693498
$ws.Range("E19").Value = 1.024235045498972
$ws.Range("F19").Value = 1.030137244205583
$ws.Range("I19").Value = 1.029452250875537
$ws.Range("J19").Value = 1.026642001382887
$ws.Range("K19").Value = 1.028771338553647
$ws.Range("L19").Value = 1.027585457575391
$ws.Range("M19").Value = 1.033467067932412
$ws.Range("N19").Value = 1.028099950782335
# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.020094100610118
$ws.Range("D20").Value = 1.025163550920295
$ws.Range("E20").Value = 1.023886272204082
$ws.Range("F20").Value = 1.029702831373378
$ws.Range("I20").Value = 1.029386212247755
$ws.Range("J20").Value = 1.026399254937919
$ws.Range("K20").Value = 1.028575389777678
$ws.Range("L20").Value = 1.027302684161751
$ws.Range("M20").Value = 1.03309851926944
$ws.Range("N20").Value = 1.027856859609571
# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.018898094577186
$ws.Range("D21").Value = 1.02431309391219
$ws.Range("E21").Value = 1.022753031440583
$ws.Range("F21").Value = 1.028290774924524
$ws.Range("I21").Value = 1.029169575711434
$ws.Range("J21").Value = 1.025609458496372
$ws.Range("K21").Value = 1.027937276035802
$ws.Range("L21").Value = 1.026383177345645
$ws.Range("M21").Value = 1.031899841990674
$ws.Range("N21").Value = 1.027065941566486
# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.018146208073582
$ws.Range("D22").Value = 1.023778272469284
$ws.Range("E22").Value = 1.022041019221768
$ws.Range("F22").Value = 1.027403151202976
$ws.Range("I22").Value = 1.029031882088619
$ws.Range("J22").Value = 1.025112417997349
$ws.Range("K22").Value = 1.027535252383724
$ws.Range("L22").Value = 1.02580490567135
$ws.Range("M22").Value = 1.031145803850266
$ws.Range("N22").Value = 1.026568195212943
# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.018544789056472
$ws.Range("D23").Value = 1.024061801302201
$ws.Range("E23").Value = 1.022418422911396
$ws.Range("F23").Value = 1.027873679653904
$ws.Range("I23").Value = 1.029105016884953
$ws.Range("J23").Value = 1.025375952625593
$ws.Range("K23").Value = 1.027748450072709
$ws.Range("L23").Value = 1.026111471873219
$ws.Range("M23").Value = 1.031545570387969
$ws.Range("N23").Value = 1.026832104090587
# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.020114153391225
$ws.Range("D24").Value = 1.025177807132468
$ws.Range("E24").Value = 1.02390527964503
$ws.Range("F24").Value = 1.029726508067079
$ws.Range("I24").Value = 1.029389818999432
$ws.Range("J24").Value = 1.026412488126422
$ws.Range("K24").Value = 1.028586073998029
$ws.Range("L24").Value = 1.027318097434783
$ws.Range("M24").Value = 1.033118608870963
$ws.Range("N24").Value = 1.02787011159072
# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.021935729908263
$ws.Range("D25").Value = 1.026472398639084
$ws.Range("E25").Value = 1.025632881640839
$ws.Range("F25").Value = 1.031877503504373
$ws.Range("I25").Value = 1.029713845066794
$ws.Range("J25").Value = 1.0276133149435
$ws.Range("K25").Value = 1.029554533380393
$ws.Range("L25").Value = 1.028717709886815
$ws.Range("M25").Value = 1.03494241132082
$ws.Range("N25").Value = 1.029072643719614
